$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.43126688834027505
$ws.Range("B1").Value = 0.43015425422417763
$ws.Range("A2").Value = -0.31669422909186196
$ws.Range("B2").Value = 0.31332109363173011
$ws.Range("A3").Value = -0.21036876431972829
$ws.Range("B3").Value = 0.2092897069877111
$ws.Range("A4").Value = -0.19728970716584193
$ws.Range("B4").Value = 0.19630152140138613
$ws.Range("A5").Value = -0.15171220446882838
$ws.Range("B5").Value = 0.14989684247922774
$ws.Range("A6").Value = -0.064761874374243966
$ws.Range("B6").Value = 0.064697713631348019
$ws.Range("A7").Value = -0.044697714407195832
$ws.Range("B7").Value = 0.044583960467697636
$ws.Range("A8").Value = -0.02458396124798945
$ws.Range("B8").Value = 0.024535228558585587
$ws.Range("A9").Value = -0.018535229231031458
$ws.Range("B9").Value = 0.018506754094679501
$ws.Range("A10").Value = -0.012506754769553652
$ws.Range("B10").Value = 0.012507988267728365
$ws.Range("A11").Value = -0.008007988930724963
$ws.Range("B11").Value = 0.0080083366109811038
$ws.Range("A12").Value = -0.0020083372859382997
$ws.Range("B12").Value = 0.0020077379974119403
$ws.Range("A13").Value = 0.0039922613277107999
$ws.Range("B13").Value = -0.0039931562630410866
$ws.Range("A14").Value = 0.015993155540865445
$ws.Range("B14").Value = -0.016008556516587724
$ws.Range("A15").Value = -0.021052610546063555
$ws.Range("B15").Value = 0.021027501958152506
$ws.Range("A16").Value = -0.015027502633362388
$ws.Range("B16").Value = 0.01500420347567788
$ws.Range("A17").Value = -0.0090042041536930739
$ws.Range("B17").Value = 0.0089999992975462462
$ws.Range("A18").Value = -0.10346002576376989
$ws.Range("B18").Value = 0.10335403845758862
$ws.Range("A19").Value = -0.094354039082727326
$ws.Range("B19").Value = 0.093508556095184492
$ws.Range("A20").Value = -0.018013490172787527
$ws.Range("B20").Value = 0.018004280969647013
$ws.Range("A21").Value = -0.009004281616835641
$ws.Range("B21").Value = 0.0089999993522216215
$ws.Range("A22").Value = -0.093952329958225889
$ws.Range("B22").Value = 0.093638005214083719
$ws.Range("A23").Value = -0.084638005866168875
$ws.Range("B23").Value = 0.084127415206445733
$ws.Range("A24").Value = -0.042127416131100937
$ws.Range("B24").Value = 0.041999999070143623
$ws.Range("A25").Value = -0.079134968767760938
$ws.Range("B25").Value = 0.078987849722814474
$ws.Range("A26").Value = -0.072987850387470132
$ws.Range("B26").Value = 0.072801855282982331
$ws.Range("A27").Value = -0.066801855950177291
$ws.Range("B27").Value = 0.066178552510988808
$ws.Range("A28").Value = -0.060178553189655482
$ws.Range("B28").Value = 0.059761769549305122
$ws.Range("A29").Value = -0.047761770283177185
$ws.Range("B29").Value = 0.04757523270203734
$ws.Range("A30").Value = -0.042171616627185493
$ws.Range("B30").Value = 0.042019388177807038
$ws.Range("A31").Value = -0.027019388945983991
$ws.Range("B31").Value = 0.027000612986180172
$ws.Range("A32").Value = -0.0060006138032129286
$ws.Range("B32").Value = 0.005999999302053638
